$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for A2:A49 (replacing the prior HKG run's results).
$newValues = @(
    87.39746457867264,
    87.695749440715886,
    88.441461595824009,
    91.349739000745714,
    91.648023862788961,
    91.648023862788961,
    78.374347501864278,
    80.536912751677846,
    80.686055182699477,
    78.523489932885909,
    81.282624906785983,
    84.041759880686044,
    79.791200596569723,
    79.716629381058908,
    82.102908277404921,
    81.058911260253538,
    83.743475018642798,
    87.322893363161825,
    88.73974645786727,
    89.112602535421331,
    88.814317673378071,
    87.173750932140194,
    90.007457121551084,
    89.932885906040269,
    84.41461595824012,
    84.190902311707688,
    84.489187173750935,
    82.102908277404921,
    81.282624906785983,
    82.177479492915737,
    88.292319164802393,
    94.034302759134974,
    91.648023862788961,
    86.577181208053688,
    93.586875466070097,
    83.22147651006712,
    90.23117076808353,
    86.204325130499626,
    87.099179716629379,
    81.208053691275168,
    81.357196122296799,
    79.716629381058908,
    80.462341536167031,
    84.787472035794181,
    84.041759880686044,
    81.133482475764353,
    81.133482475764353,
    81.357196122296799
)

# Re-assert the header/label in A1 (shared string "HK_R_acc_LT").
$ws.Range("A1").Value = "HK_R_acc_LT"

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
